$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (2022Q4)
$ws.Range("B9").Value = 0.0746880969388026
$ws.Range("C9").Value = 0.3688458615154391
$ws.Range("D9").Value = 5.763500000000001
$ws.Range("E9").Value = -0.2941577645766364
$ws.Range("F9").Value = 0.2582411831116511

# Row 10 (2023Q1)
$ws.Range("B10").Value = 0.07712709796616259
$ws.Range("C10").Value = 0.3622564205410629
$ws.Range("D10").Value = 4.347
$ws.Range("E10").Value = -0.2851293225749004
$ws.Range("F10").Value = -0.04885964987984415

# Row 13 (2023Q4)
$ws.Range("F13").Value = -0.03601842022030477

# Row 14 (2024Q1)
$ws.Range("F14").Value = -0.02941602417012457
